$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 665
$ws1.Range("F4").Value = 1554
$ws1.Range("F5").Value = 718
$ws1.Range("F6").Value = 25

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 665
$ws4.Range("F4").Value = 1554
$ws4.Range("F6").Value = 718
$ws4.Range("F7").Value = 25
